# Generate Report for Handback
#
# Re-runs the handback-status report generation: the first entry's GUID
# becomes af01bf42-34aa-4cbf-876c-d360ca1fd7b9, the second entry's GUID
# becomes fffffebc3136-fb07-454e-9d70-6fa9d02057a0, the handoff/handback
# content-hash is refreshed to 104ac376faed4118def4e901bf122f240ac8750d,
# and the handoff/handback timestamps move a minute or so later.

$wb = $excel.ActiveWorkbook

$newMd1 = "af01bf42-34aa-4cbf-876c-d360ca1fd7b9.md"
$newMd2 = "fffffebc3136-fb07-454e-9d70-6fa9d02057a0.md"

$newXlfZh = "af01bf42-34aa-4cbf-876c-d360ca1fd7b9.104ac376faed4118def4e901bf122f240ac8750d.zh-cn.xlf"
$newXlfDe = "af01bf42-34aa-4cbf-876c-d360ca1fd7b9.104ac376faed4118def4e901bf122f240ac8750d.de-de.xlf"

$newHandoffZh = "2016-03-13 00:55:10"
$newHandbackZh = "2016-03-13 00:55:26"
$newHandoffDe = "2016-03-13 00:55:14"
$newHandbackDe = "2016-03-13 00:55:32"

# ---------------------------------------------------------------------
# Sheet "Overview": File Name column holds the two .md entries.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMd1
$wsOverview.Range("A3").Value = $newMd2

foreach ($h in $wsOverview.Hyperlinks) {
    if ($h.Range.Row -eq 2 -and $h.Range.Column -eq 1) {
        $h.TextToDisplay = $newMd1
    }
    elseif ($h.Range.Row -eq 3 -and $h.Range.Column -eq 1) {
        $h.TextToDisplay = $newMd2
    }
}

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newMd1
$wsZh.Range("D2").Value = $newXlfZh
$wsZh.Range("E2").Value = $newHandoffZh
$wsZh.Range("F2").Value = $newMd1
$wsZh.Range("G2").Value = $newXlfZh
$wsZh.Range("H2").Value = $newHandbackZh

$wsZh.Range("A3").Value = $newMd2
$wsZh.Range("D3").Value = $newXlfZh
$wsZh.Range("E3").Value = $newHandoffZh
$wsZh.Range("F3").Value = $newMd2
$wsZh.Range("G3").Value = $newXlfZh
$wsZh.Range("H3").Value = $newHandbackZh

foreach ($h in $wsZh.Hyperlinks) {
    $r = $h.Range.Row
    $c = $h.Range.Column
    if ($c -eq 1 -and $r -eq 2) { $h.TextToDisplay = $newMd1 }
    elseif ($c -eq 1 -and $r -eq 3) { $h.TextToDisplay = $newMd2 }
    elseif ($c -eq 4 -and ($r -eq 2 -or $r -eq 3)) { $h.TextToDisplay = $newXlfZh }
    elseif ($c -eq 6 -and $r -eq 2) { $h.TextToDisplay = $newMd1 }
    elseif ($c -eq 6 -and $r -eq 3) { $h.TextToDisplay = $newMd2 }
    elseif ($c -eq 7 -and ($r -eq 2 -or $r -eq 3)) { $h.TextToDisplay = $newXlfZh }
}

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newMd1
$wsDe.Range("D2").Value = $newXlfDe
$wsDe.Range("E2").Value = $newHandoffDe
$wsDe.Range("F2").Value = $newMd1
$wsDe.Range("G2").Value = $newXlfDe
$wsDe.Range("H2").Value = $newHandbackDe

$wsDe.Range("A3").Value = $newMd2
$wsDe.Range("D3").Value = $newXlfDe
$wsDe.Range("E3").Value = $newHandoffDe
$wsDe.Range("F3").Value = $newMd2
$wsDe.Range("G3").Value = $newXlfDe
$wsDe.Range("H3").Value = $newHandbackDe

foreach ($h in $wsDe.Hyperlinks) {
    $r = $h.Range.Row
    $c = $h.Range.Column
    if ($c -eq 1 -and $r -eq 2) { $h.TextToDisplay = $newMd1 }
    elseif ($c -eq 1 -and $r -eq 3) { $h.TextToDisplay = $newMd2 }
    elseif ($c -eq 4 -and ($r -eq 2 -or $r -eq 3)) { $h.TextToDisplay = $newXlfDe }
    elseif ($c -eq 6 -and $r -eq 2) { $h.TextToDisplay = $newMd1 }
    elseif ($c -eq 6 -and $r -eq 3) { $h.TextToDisplay = $newMd2 }
    elseif ($c -eq 7 -and ($r -eq 2 -or $r -eq 3)) { $h.TextToDisplay = $newXlfDe }
}
